# ---------------------------------------------------------------------------
# Adds a new "CUDA, only for dTAd" results block (rows 151-169, plus some
# trailing formatted filler rows 171-186) to the existing benchmark sheet,
# mirroring the layout already used for the "OPENMP" block (rows 114-132).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Shared strings must be interned in this exact order so that they line
#    up with the indices used throughout the new block:
#       14 -> "CUDA，仅dtad"
#       15 -> "cuda"
#       16 -> "cuda总"
#       17 -> "cuda部分"
# ---------------------------------------------------------------------------
$ws.Range("A151").Value = "CUDA，仅dtad"
$ws.Range("H153").Value = "cuda"
$ws.Range("C153").Value = "cuda总"
$ws.Range("D153").Value = "cuda部分"

# reuse the strings just created for their remaining occurrences
$ws.Range("I153").Value = "cuda部分"
$ws.Range("M153").Value = "cuda"
$ws.Range("N153").Value = "cuda部分"

# ---------------------------------------------------------------------------
# 2. Section title row (151), formatted like the "OPENMP" title (A114/B114)
#    but centered, and merged across A151:C151.
# ---------------------------------------------------------------------------
$ws.Range("B114").Copy()
$ws.Range("A151:C151").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A151").Value = "CUDA，仅dtad"
$ws.Range("A151").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A151:C151").Merge()

# ---------------------------------------------------------------------------
# 3. Column-group label row (152): mirrors row 115 / row 78 / row 1 layout -
#    "运行时间，单位(s)" / "迭代次数" / "平均每次迭代时间，单位(ms)"
# ---------------------------------------------------------------------------
$ws.Range("B115:D115").Copy()
$ws.Range("B152:D152").PasteSpecial(-4122)
$ws.Range("G115:I115").Copy()
$ws.Range("G152:I152").PasteSpecial(-4122)
$ws.Range("L115:N115").Copy()
$ws.Range("L152:N152").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B152").Value = "运行时间，单位(s)"
$ws.Range("G152").Value = "迭代次数"
$ws.Range("L152").Value = "平均每次迭代时间，单位(ms)"

$ws.Range("B152:D152").Merge()
$ws.Range("G152:I152").Merge()
$ws.Range("L152:N152").Merge()

# ---------------------------------------------------------------------------
# 4. Detail header row (153): "问题规模" / "base" / "cuda总" / "cuda部分" etc,
#    formatting copied from the analogous row 116 cells. The "cuda" value
#    columns (C,D / H,I / M,N) additionally pick up the row-114 title font.
# ---------------------------------------------------------------------------
$ws.Range("A116:E116").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("F116:J116").Copy()
$ws.Range("F153:J153").PasteSpecial(-4122)
$ws.Range("K116:O116").Copy()
$ws.Range("K153:O153").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B114").Copy()
$ws.Range("C153:D153").PasteSpecial(-4122)
$ws.Range("H153:I153").PasteSpecial(-4122)
$ws.Range("M153:N153").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A153").Value = "问题规模"
$ws.Range("B153").Value = "base"
$ws.Range("C153").Value = "cuda总"
$ws.Range("D153").Value = "cuda部分"
$ws.Range("F153").Value = "问题规模"
$ws.Range("G153").Value = "base"
$ws.Range("H153").Value = "cuda"
$ws.Range("I153").Value = "cuda部分"
$ws.Range("K153").Value = "问题规模"
$ws.Range("L153").Value = "base"
$ws.Range("M153").Value = "cuda"
$ws.Range("N153").Value = "cuda部分"

# ---------------------------------------------------------------------------
# 5. Data rows 154-169 (problem sizes 128..2048), three data groups:
#      A:C  -> size, dTAd total time, cuda total time
#      F:H  -> size, dTAd iteration count, cuda iteration count
#      K:M  -> size, per-iteration time (ms) formulas
# ---------------------------------------------------------------------------

$sizes  = 128,256,384,512,640,768,896,1024,1152,1280,1408,1536,1664,1792,1920,2048
$bVals  = 0.14218838,1.36119203,4.1844674407301401,12.182143569999999,25.184467440733901,31.868529219999999,55.184467440734601,138.184467440731,99.705498199999994,212.67718826999999,232.18446744073501,419.18446744073401,429.184467440736,638.14885233999996,994.18446744073503,1161.4684852
$cVals  = 0.70580854999999998,2.9186605999999999,3.6004505600000001,5.5624348000000001,15.1844674407313,23.498286581999999,43.184467440734203,51.184467440733897,87.560694255000001,70.197168245,140.4751784,136.63108869999999,190.18446744073401,147.44270319200001,355.18446744073202,422.53599036000003
$gVals  = 676,1659,1825,3740,5097,4414,5802,10906,6345,10818,9693,14670,13141,16652,22916,23799
$hVals  = 616,2115,2789,2979,5323,6751,8417,8708,11567,7942,12640,10706,12805,8719,18705,19499

for ($i = 0; $i -lt 16; $i++) {
    $r = 154 + $i
    $ws.Range("A$r").Value = $sizes[$i]
    $ws.Range("B$r").Value = $bVals[$i]
    $ws.Range("C$r").Value = $cVals[$i]
    $ws.Range("F$r").Value = $sizes[$i]
    $ws.Range("G$r").Value = $gVals[$i]
    $ws.Range("H$r").Value = $hVals[$i]
    $ws.Range("K$r").Value = $sizes[$i]
    $ws.Range("L$r").Formula = "=B$r/G$r*1000"
    $ws.Range("M$r").Formula = "=C$r/H$r*1000"
}

# ---------------------------------------------------------------------------
# 6. Trailing formatted filler cells C171:C186 (elapsed-time number format,
#    same style already used e.g. at Q132 / E99 and friends).
# ---------------------------------------------------------------------------
$ws.Range("Q132").Copy()
$ws.Range("C171:C186").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 7. Refresh the window selection to match the final saved view.
# ---------------------------------------------------------------------------
$ws.Range("O157").Select()
